# The workbook has "Sheet1" (UserID list, currently 1 data row) and "Sheet2"
# (same columns, 2 data rows). This edit appends Sheet2's two data rows to the
# bottom of Sheet1, i.e. it copies Sheet2!A1:G2 into Sheet1 starting at A3,
# carrying over values, number/cell formatting (so the "Email" column keeps
# its hyperlink-looking style), and then restores the live hyperlink that
# only existed on the second of those two source rows (Sheet2!E2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sheet1 (destination)
$ws2 = $wb.Worksheets.Item(2)   # Sheet2 (source)

# Copy values + formatting for the 2 rows of Sheet2 into Sheet1, right below
# the existing data (new rows 3 and 4).
$ws2.Range("A1:G2").Copy($ws1.Range("A3"))

# Only Sheet2's second row (now Sheet1 row 4) had a real hyperlink on its
# Email cell (mailto link), so recreate that relationship on E4 and make
# sure the cell keeps using the built-in "Hyperlink" style afterward.
$ws1.Hyperlinks.Add($ws1.Range("E4"), "mailto:etdg@gmail.com")
$ws1.Range("E4").Style = "Hyperlink"

# Match the selection Excel leaves behind after pasting the new rows.
$ws1.Range("A3:G4").Select()
